$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value, derived from the cryptos-list refresh diff.
$updates = [ordered]@{
    "D2" = "67.955.42"
    "E2" = "  +1.42%  "
    "D3" = "3.556.71"
    "E3" = "  -0.14%  "
    "E4" = "  +0.19%  "
    "D5" = "205.09"
    "E5" = "  +7.86%  "
    "D6" = "557.77"
    "E6" = "  -1.78%  "
    "D7" = "0.604"
    "E7" = "  -1.48%  "
    "E8" = "  -0.06%  "
    "D9" = "0.668"
    "D10" = "62.44"
    "E10" = "  +12.05%  "
    "D11" = "0.145"
    "E11" = "  -2.64%  "
    "D12" = "0.0000275"
    "E12" = "  +2.41%  "
    "D13" = "9.99"
    "E13" = "  +1.68%  "
    "D14" = "4.142.94"
    "E14" = "  +0.00%  "
    "D15" = "3.578.20"
    "E15" = "  +0.20%  "
    "E16" = "  -0.15%  "
    "D17" = "18.87"
    "E17" = "  +4.16%  "
    "D18" = "67.827.39"
    "D19" = "12.02"
    "E19" = "  -0.93%  "
    "D20" = "1.05"
    "E20" = "  -0.77%  "
    "D21" = "392.47"
    "E21" = "  -1.87%  "
    "D22" = "4.08"
    "E22" = "  -1.44%  "
    "D23" = "12.13"
    "E23" = "  +2.49%  "
    "D24" = "83.66"
    "E24" = "  -2.23%  "
    "B25" = "InternetComputer(DFINITY)"
    "C25" = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
    "D25" = "12.27"
    "E25" = "  -1.39%  "
    "B26" = "ImmutableX"
    "C26" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
    "D26" = "2.83"
    "E26" = "  -2.31%  "
    "D27" = "3.81"
    "E27" = "  +4.11%  "
    "D28" = "8.98"
    "E28" = "  +0.50%  "
    "D29" = "717.73"
    "E29" = "  +12.51%  "
    "D30" = "31.04"
    "E30" = "  -0.13%  "
    "D31" = "7.45"
    "E31" = "  -4.43%  "
    "D32" = "11.91"
    "E32" = "  -1.43%  "
    "D33" = "63.59"
    "E33" = "  -0.14%  "
    "E34" = "  -2.60%  "
    "D35" = "40.62"
    "E35" = "  -3.80%  "
    "D36" = "0.414"
    "E36" = "  +2.67%  "
    "D37" = "1.00"
    "E37" = "  +0.25%  "
    "D38" = "3.19"
    "E38" = "  +6.74%  "
    "B39" = "Maker"
    "C39" = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
    "D39" = "3.137.84"
    "E39" = "  -1.50%  "
    "B40" = "dogwifhat"
    "C40" = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
    "D40" = "3.05"
    "E40" = "  +26.00%  "
    "D41" = "0.0₃0724"
    "E41" = "  -4.63%  "
    "B42" = "FirstDigitalUSD"
    "C42" = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
    "D42" = "1.00"
    "E42" = "  -0.18%  "
    "B43" = "Kaspa"
    "C43" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "D43" = "0.130"
    "E43" = "  -1.92%  "
    "D44" = "2.53"
    "E44" = "  -6.41%  "
    "D45" = "2.74"
    "E45" = "  +9.02%  "
    "D46" = "0.0406"
    "E46" = "  -1.47%  "
    "B47" = "Stellar"
    "C47" = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
    "D47" = "0.129"
    "E47" = "  -0.65%  "
    "B48" = "ApeXProtocol"
    "C48" = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
    "D48" = "3.03"
    "E48" = "  -3.03%  "
    "B49" = "THORChain"
    "C49" = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
    "D49" = "8.52"
    "E49" = "  +0.50%  "
    "D50" = "137.93"
    "E50" = "  -2.49%  "
    "D51" = "2.66"
    "E51" = "  -1.16%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (prices like "205.09")
    # and percent strings keep their exact characters instead of becoming
    # General-formatted numbers, then drop back to the default (unstyled)
    # format so no stray style index is left on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
